# Updates cryptos list prices / 1h volume % (and a few re-ranked coin rows)
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
# Values that look like plain numbers are prefixed with a leading apostrophe so
# Excel stores them as text (matching the source data, e.g. "606.81"), just like
# typing '606.81 directly into a cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.325.58"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.605.63"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'606.81"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'140.41"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("D7").Value = "3.602.96"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  +5.60%  "
$ws.Range("E12").Value = "  +2.51%  "
$ws.Range("D13").Value = "4.224.00"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "'28.76"
$ws.Range("E14").Value = "  +5.72%  "
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "3.608.10"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.117"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "66.360.10"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "'10.23"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'14.77"
$ws.Range("E20").Value = "  +3.45%  "
$ws.Range("D21").Value = "'5.95"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'399.77"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "'0.595"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").Value = "3.752.39"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "'75.22"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +5.31%  "
$ws.Range("D28").Value = "'8.21"
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("D29").Value = "'1.66"
$ws.Range("E29").Value = "  +28.23%  "
$ws.Range("E30").Value = "  +7.21%  "
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").Value = "'0.997"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "3.612.93"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "'24.77"
$ws.Range("E34").Value = "  +3.86%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.66"
$ws.Range("E37").Value = "  +5.91%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.45"
$ws.Range("E38").Value = "  +9.09%  "
$ws.Range("D39").Value = "'7.13"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").Value = "'168.91"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("E41").Value = "  +5.40%  "
$ws.Range("D42").Value = "'0.846"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").Value = "'26.79"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("E44").Value = "  +7.76%  "
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").Value = "'43.19"
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.74"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").Value = "'7.08"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("D50").Value = "2.470.05"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'22.96"
$ws.Range("E51").Value = "  +9.94%  "
